$d = $word.ActiveDocument

# Helper: insert a literal <w:p>...</w:p> fragment by replacing the given
# (currently-empty) placeholder paragraph range with it via raw OOXML.
# Word's InsertXML on a paragraph range inserts the new paragraph *before*
# the placeholder and leaves the (now empty) placeholder paragraph behind,
# so callers keep using the newly-created trailing placeholder for the next
# insertion and drop it once every paragraph has been inserted.
function Insert-RawParagraph($placeholderRange, [string]$innerXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $placeholderRange.InsertXML($pkg) | Out-Null
}

# Anchor: the last paragraph in the document ("Api/deleteEmployee etc").
$anchor = $d.Paragraphs.Last
$anchor.Range.InsertParagraphAfter() | Out-Null
$placeholder = $d.Paragraphs.Last

# --- "Spring Data JPA" (plain paragraph, no list) ---
$p1Xml = '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Spring Data JPA</w:t></w:r></w:p>'
Insert-RawParagraph $placeholder.Range $p1Xml
$placeholder = $d.Paragraphs.Last

# --- "Helps reduce code amount" (ListParagraph, ilvl 0, numId 1) ---
$p2Xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Helps reduce code amount</w:t></w:r></w:p>'
Insert-RawParagraph $placeholder.Range $p2Xml
$placeholder = $d.Paragraphs.Last

# --- "Plug in the specific DAO and spring will supply the CRUD implementation for you" ---
$p3Xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Plug in the specific DAO and spring will supply the CRUD implementation for you</w:t></w:r></w:p>'
Insert-RawParagraph $placeholder.Range $p3Xml
$placeholder = $d.Paragraphs.Last

# --- "Use JPARepository interface for these features/operations" ---
$p4Xml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Use JPARepository interface for these features/operations</w:t></w:r></w:p>'
Insert-RawParagraph $placeholder.Range $p4Xml

# Drop the trailing empty placeholder paragraph left over from the chain.
# (A plain zero-length Range.Delete() on the very last paragraph in the
# document is a no-op in this runtime, so instead extend the deletion
# range to also swallow the preceding paragraph mark, which reliably
# removes the empty trailing paragraph.)
$count = $d.Paragraphs.Count
$trailing = $d.Paragraphs.Item($count)
$prev = $d.Paragraphs.Item($count - 1)
$cleanupRange = $d.Range($prev.Range.End - 1, $trailing.Range.End)
$cleanupRange.Delete() | Out-Null

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
